# monster.xlsx — Sheet1 header row tweak + selection move.
#
# Before: B1 = "x " (trailing space), C1 = "y"
# After : B1 = "x"                 , C1 = "y"
# (the shared-string table ends up de-duplicated/reordered so that the
# stray trailing-space variant of "x" is gone and "y"/"x" occupy the
# remaining two slots — re-setting the two header cells to their clean
# text values reproduces that.)
#
# The saved selection/active cell also moves from D6 to H10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "x"
$ws.Range("C1").Value = "y"

$ws.Range("H10").Select()
